# Comfenalco Cartagena - Estado de Cuenta
# "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# Adds a new payroll-debt period (2509) as a new data row, updates the
# running totals (VALOR MORA / Cant. Periodos) accordingly, and keeps the
# two closing signature lines directly below the (now one row taller)
# data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Make room for the new period row. The data table currently ends
#    at row 21 (period 2508); push everything from row 22 down
#    (the two signature/footer rows) down by one row.
# ------------------------------------------------------------------

# Capture the current last data row's formatting (it carries the
# table's closing/bottom border) before we disturb anything, so we can
# stamp it onto the new last row.
$ws.Range("B21:J21").Copy()
$ws.Rows.Item(22).Insert()
$ws.Range("B22:J22").PasteSpecial(-4122)

# Row 21 is no longer the last row of the table, so it goes back to
# the plain "interior row" look (same as rows 16-20).
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2. Fill in the new period (2509) data row.
# ------------------------------------------------------------------
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "45447483"
$ws.Range("D22").Value = "NANCY ELVIRA MACHACON TORRES"
$ws.Range("E22").Value = "2509"
$ws.Range("F22").Value = 55042
$ws.Range("G22").Value = 1423500

# ------------------------------------------------------------------
# 3. The "Periodo Mora" column reads clearer centered - apply that to
#    the whole column of period codes, old rows and the new one alike.
# ------------------------------------------------------------------
$ws.Range("E16:E22").HorizontalAlignment = -4108

# ------------------------------------------------------------------
# 4. Refresh the summary figures at the top of the sheet.
# ------------------------------------------------------------------
$ws.Range("F13").Value = 7
$ws.Range("E11").Value = 396682

# ------------------------------------------------------------------
# 5. Keep selection tidy (matches the template's usual default).
# ------------------------------------------------------------------
$ws.Range("D2").Select()
